$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'30.408.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.01%  "

# Row 3
$ws.Range("D3").Value = "'2.127.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.57%  "

# Row 4
$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.75%  "

# Row 5
$ws.Range("D5").Value = "'347.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.95%  "

# Row 6
$ws.Range("D6").Value = "'1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.57%  "

# Row 7
$ws.Range("D7").Value = "'0.5219"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.91%  "

# Row 8
$ws.Range("D8").Value = "'0.4489"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.57%  "

# Row 9
$ws.Range("D9").Value = "'53.86"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.54%  "

# Row 10
$ws.Range("D10").Value = "'0.09393"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.04%  "

# Row 11
$ws.Range("D11").Value = "'1.188"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.54%  "

# Row 12
$ws.Range("D12").Value = "'25.50"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.93%  "

# Row 13
$ws.Range("D13").Value = "'8.702"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.26%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "'2.112.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.38%  "

# Row 15
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'6.981"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.88%  "

# Row 16
$ws.Range("D16").Value = "'103.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.02%  "

# Row 17
$ws.Range("E17").Value = "  +1.31%  "

# Row 18
$ws.Range("D18").Value = "'1.005"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.64%  "

# Row 19
$ws.Range("D19").Value = "'21.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.51%  "

# Row 20
$ws.Range("D20").Value = "'0.06703"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.10%  "

# Row 21
$ws.Range("D21").Value = "'6.323"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.46%  "

# Row 22
$ws.Range("D22").Value = "'1.005"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.41%  "

# Row 23
$ws.Range("D23").Value = "'30.374.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.49%  "

# Row 24
$ws.Range("D24").Value = "'12.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.16%  "

# Row 25
$ws.Range("D25").Value = "'2.330"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.27%  "

# Row 26
$ws.Range("D26").Value = "'2.375.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.04%  "

# Row 27
$ws.Range("E27").Value = "  +1.21%  "

# Row 28
$ws.Range("E28").Value = "  +1.61%  "

# Row 29
$ws.Range("D29").Value = "'163.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.22%  "

# Row 30
$ws.Range("D30").Value = "'134.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.74%  "

# Row 31
$ws.Range("D31").Value = "'1.168"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.32%  "

# Row 32
$ws.Range("D32").Value = "'1.815"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +12.71%  "

# Row 33
$ws.Range("D33").Value = "'0.1060"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.40%  "

# Row 34
$ws.Range("D34").Value = "'6.858"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.69%  "

# Row 35
$ws.Range("D35").Value = "'6.322"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.80%  "

# Row 36
$ws.Range("D36").Value = "'3.959"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.07%  "

# Row 37
$ws.Range("D37").Value = "'10.75"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.59%  "

# Row 38
$ws.Range("D38").Value = "'0.02648"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.43%  "

# Row 39
$ws.Range("E39").Value = "  +2.25%  "

# Row 40
$ws.Range("D40").Value = "'0.7175"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.35%  "

# Row 41
$ws.Range("E41").Value = "  +3.22%  "

# Row 42
$ws.Range("D42").Value = "'0.2260"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.49%  "

# Row 43
$ws.Range("D43").Value = "'1.341"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.73%  "

# Row 44
$ws.Range("D44").Value = "'0.6984"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.25%  "

# Row 45
$ws.Range("E45").Value = "  +4.81%  "

# Row 46
$ws.Range("D46").Value = "'2.409"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.00%  "

# Row 47
$ws.Range("D47").Value = "'1.004"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.61%  "

# Row 48
$ws.Range("E48").Value = "  +8.50%  "

# Row 49
$ws.Range("D49").Value = "'3.634"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.16%  "

# Row 50
$ws.Range("E50").Value = "  +2.90%  "

# Row 51
$ws.Range("D51").Value = "'1.231"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.25%  "

